$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target = 50 block (rows 23-27)
$ws.Range("B23").Formula = "=235.2-50"
$ws.Range("C23").Formula = "=236.4-50"
$ws.Range("D23").Formula = "=234.3-50"

$ws.Range("B24").Formula = "=248.9-50"
$ws.Range("C24").Formula = "=249.1-50"
$ws.Range("D24").Formula = "=251.6-50"

$ws.Range("B25").Formula = "=195.1-50"
$ws.Range("C25").Formula = "=196.8-50"
$ws.Range("D25").Formula = "=196.8-50"

$ws.Range("B26").Formula = "=146.3-50"
$ws.Range("C26").Formula = "=146.7-50"
$ws.Range("D26").Formula = "=147.4-50"

$ws.Range("B27").Formula = "=97.1-50"
$ws.Range("C27").Formula = "=96-50"
$ws.Range("D27").Formula = "=96.5-50"

# Target = 60 block (rows 30-34)
$ws.Range("B30").Formula = "=248.6-50"
$ws.Range("C30").Formula = "=250.1-50"
$ws.Range("D30").Formula = "=249.9-50"

$ws.Range("B31").Formula = "=247-50"
$ws.Range("C31").Formula = "=247.8-50"
$ws.Range("D31").Formula = "=248.1-50"

$ws.Range("B32").Formula = "=196.4-50"
$ws.Range("C32").Formula = "=195.9-50"
$ws.Range("D32").Formula = "=194.9-50"

$ws.Range("B33").Formula = "=145.8-50"
$ws.Range("C33").Formula = "=146.7-50"
$ws.Range("D33").Formula = "=146.1-50"

$ws.Range("B34").Formula = "=96.3-50"
$ws.Range("C34").Formula = "=96.4-50"
$ws.Range("D34").Formula = "=96.2-50"

# Target = 70 block (rows 37-41)
$ws.Range("B37").Formula = "=255.3-50"
$ws.Range("C37").Formula = "=255.8-50"
$ws.Range("D37").Formula = "=255.9-50"

$ws.Range("B38").Formula = "=251.9-50"
$ws.Range("C38").Formula = "=252.2-50"
$ws.Range("D38").Formula = "=250.5-50"

$ws.Range("B39").Formula = "=199.9-50"
$ws.Range("C39").Formula = "=200-50"
$ws.Range("D39").Formula = "=201.3-50"

$ws.Range("B40").Formula = "=149.6-50"
$ws.Range("C40").Formula = "=149.9-50"
$ws.Range("D40").Formula = "=150.1-50"

$ws.Range("B41").Formula = "=97.6-50"
$ws.Range("C41").Formula = "=98.4-50"
$ws.Range("D41").Formula = "=97.6-50"

# Move selection/active cell to reflect where the user ended up working
$ws.Range("F39").Select()
